$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rich-text edits: replace the changed substring in place ---
# A8 shared string: "Volume 31   Number  39" -> "...  40"
$ws.Range("A8").Characters(21, 2).Text = "40"

# C9 shared string: "Report Covering the Week  9/23/2024  Through  9/29/2024"
#                -> "Report Covering the Week  9/30/2024  Through  10/6/2024"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "9/30/2024"
$c9.Characters(47, 9).Text = "10/6/2024"

# --- Cells whose data type flips between number and text placeholder ---
# ($ws.Range(template).Copy(...) clones the correct style (incl. number format),
#  then the real value/text is written on top of it.)
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 4
$ws.Range("D14").Copy($ws.Range("C20"))
$ws.Range("C20").Formula = "'0"
$ws.Range("D14").Copy($ws.Range("D20"))
$ws.Range("D20").Formula = "'0"
$ws.Range("D14").Copy($ws.Range("E20"))
$ws.Range("E20").Formula = "'***.*"
$ws.Range("D14").Copy($ws.Range("C26"))
$ws.Range("C26").Formula = "'0"
$ws.Range("I14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 100
$ws.Range("I14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("I14").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("K14").Copy($ws.Range("H29"))
$ws.Range("H29").Value = -100
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("I14").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K14").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100
$ws.Range("I14").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1
$ws.Range("D14").Copy($ws.Range("F33"))
$ws.Range("F33").Formula = "'0"

# --- Plain numeric updates (style/type unchanged) ---
$ws.Range("M16").Value = -9.523809523809
$ws.Range("N16").Value = -58.695652173913
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 1.388888888888
$ws.Range("L17").Value = 97.297297297297
$ws.Range("M17").Value = 102.777777777778
$ws.Range("N17").Value = -13.095238095238
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 30
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = 7.142857142857
$ws.Range("M18").Value = -63.414634146341
$ws.Range("N18").Value = -87.654320987654
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -39.130434782608
$ws.Range("I19").Value = 217
$ws.Range("J19").Value = 219
$ws.Range("K19").Value = -0.913242009132
$ws.Range("L19").Value = 3.827751196172
$ws.Range("M19").Value = 99.082568807339
$ws.Range("N19").Value = 21.229050279329
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -33.333333333333
$ws.Range("L20").Value = -61.627906976744
$ws.Range("M20").Value = 10
$ws.Range("N20").Value = -94.043321299639
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 37.5
$ws.Range("F21").Value = 33
$ws.Range("G21").Value = 41
$ws.Range("H21").Value = -19.512195121951
$ws.Range("I21").Value = 378
$ws.Range("J21").Value = 413
$ws.Range("K21").Value = -8.474576271186
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 35
$ws.Range("N21").Value = -65.945945945946
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -66.666666666666
$ws.Range("G24").Value = 39
$ws.Range("H24").Value = -30.769230769230
$ws.Range("I24").Value = 277
$ws.Range("J24").Value = 368
$ws.Range("K24").Value = -24.728260869565
$ws.Range("L24").Value = -23.055555555555
$ws.Range("M24").Value = -34.047619047619
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 124
$ws.Range("J25").Value = 165
$ws.Range("K25").Value = -24.848484848484
$ws.Range("L25").Value = 42.528735632183
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 30.769230769230
$ws.Range("J26").Value = 154
$ws.Range("K26").Value = -17.532467532467
$ws.Range("L26").Value = -2.307692307692
$ws.Range("M26").Value = -25.730994152046
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 125
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = -18.75
$ws.Range("L28").Value = -18.75
$ws.Range("J29").Value = 2
$ws.Range("J30").Value = 2
$ws.Range("I31").Value = 3
$ws.Range("L31").Value = -40
